$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    # Preserve original cell style while forcing the assigned value to be
    # stored as literal text (avoids Excel auto-coercing numeric-looking
    # strings like "1.00" or "17.99" into numbers and losing formatting).
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '66.002.71'
Set-TextValue $ws.Range("E2") '  -0.99%  '
Set-TextValue $ws.Range("D3") '3.289.89'
Set-TextValue $ws.Range("E3") '  -2.13%  '
Set-TextValue $ws.Range("E4") '  +0.10%  '
Set-TextValue $ws.Range("D5") '187.67'
Set-TextValue $ws.Range("E5") '  +2.21%  '
Set-TextValue $ws.Range("D6") '554.42'
Set-TextValue $ws.Range("E6") '  -0.49%  '
Set-TextValue $ws.Range("D7") '1.00'
Set-TextValue $ws.Range("E7") '  +0.12%  '
Set-TextValue $ws.Range("D8") '0.584'
Set-TextValue $ws.Range("E8") '  -2.40%  '
Set-TextValue $ws.Range("D9") '3.280.19'
Set-TextValue $ws.Range("E9") '  -2.07%  '
Set-TextValue $ws.Range("E10") '  -2.12%  '
Set-TextValue $ws.Range("D11") '0.585'
Set-TextValue $ws.Range("E11") '  -1.40%  '
Set-TextValue $ws.Range("D12") '47.29'
Set-TextValue $ws.Range("E12") '  -1.24%  '
Set-TextValue $ws.Range("D13") '0.0000269'
Set-TextValue $ws.Range("E13") '  +0.40%  '
Set-TextValue $ws.Range("D14") '8.61'
Set-TextValue $ws.Range("E14") '  -1.17%  '
Set-TextValue $ws.Range("D15") '3.827.50'
Set-TextValue $ws.Range("E15") '  -1.62%  '
Set-TextValue $ws.Range("D16") '614.06'
Set-TextValue $ws.Range("E16") '  +1.66%  '
Set-TextValue $ws.Range("B17") 'WrappedBTC'
Set-TextValue $ws.Range("C17") 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range("D17") '65.981.87'
Set-TextValue $ws.Range("E17") '  -0.72%  '
Set-TextValue $ws.Range("B18") 'Chainlink'
Set-TextValue $ws.Range("C18") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D18") '17.99'
Set-TextValue $ws.Range("E18") '  +0.74%  '
Set-TextValue $ws.Range("E19") '  -0.15%  '
Set-TextValue $ws.Range("D20") '3.300.77'
Set-TextValue $ws.Range("E20") '  -1.56%  '
Set-TextValue $ws.Range("D21") '10.93'
Set-TextValue $ws.Range("E21") '  -6.23%  '
Set-TextValue $ws.Range("D22") '0.906'
Set-TextValue $ws.Range("E22") '  -0.59%  '
Set-TextValue $ws.Range("D23") '18.39'
Set-TextValue $ws.Range("E23") '  +8.93%  '
Set-TextValue $ws.Range("D24") '101.91'
Set-TextValue $ws.Range("E24") '  +4.15%  '
Set-TextValue $ws.Range("D25") '4.94'
Set-TextValue $ws.Range("E25") '  -2.62%  '
Set-TextValue $ws.Range("D26") '3.92'
Set-TextValue $ws.Range("E26") '  -3.49%  '
Set-TextValue $ws.Range("E27") '  +0.14%  '
Set-TextValue $ws.Range("E28") '  -0.73%  '
Set-TextValue $ws.Range("D29") '9.56'
Set-TextValue $ws.Range("E29") '  +1.32%  '
Set-TextValue $ws.Range("D30") '8.62'
Set-TextValue $ws.Range("E30") '  -2.32%  '
Set-TextValue $ws.Range("D31") '30.14'
Set-TextValue $ws.Range("E31") '  -2.07%  '
Set-TextValue $ws.Range("D32") '4.03'
Set-TextValue $ws.Range("E32") '  +4.87%  '
Set-TextValue $ws.Range("D33") '6.48'
Set-TextValue $ws.Range("E33") '  +2.13%  '
Set-TextValue $ws.Range("B34") 'Cosmos'
Set-TextValue $ws.Range("C34") 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D34") '11.04'
Set-TextValue $ws.Range("E34") '  -1.21%  '
Set-TextValue $ws.Range("B35") 'Bittensor'
Set-TextValue $ws.Range("C35") 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D35") '555.82'
Set-TextValue $ws.Range("E35") '  +5.74%  '
Set-TextValue $ws.Range("E36") '  -0.51%  '
Set-TextValue $ws.Range("D37") '3.818.42'
Set-TextValue $ws.Range("E37") '  -0.39%  '
Set-TextValue $ws.Range("D38") '57.23'
Set-TextValue $ws.Range("E38") '  -1.22%  '
Set-TextValue $ws.Range("E39") '  +0.12%  '
Set-TextValue $ws.Range("D40") '0.0₃0721'
Set-TextValue $ws.Range("E40") '  -0.60%  '
Set-TextValue $ws.Range("E41") '  -2.96%  '
Set-TextValue $ws.Range("D42") '33.86'
Set-TextValue $ws.Range("E42") '  +4.06%  '
Set-TextValue $ws.Range("D43") '2.70'
Set-TextValue $ws.Range("E43") '  -0.33%  '
Set-TextValue $ws.Range("E44") '  +0.80%  '
Set-TextValue $ws.Range("D45") '0.335'
Set-TextValue $ws.Range("E45") '  -4.13%  '
Set-TextValue $ws.Range("D46") '0.0419'
Set-TextValue $ws.Range("E46") '  +0.19%  '
Set-TextValue $ws.Range("D47") '3.18'
Set-TextValue $ws.Range("E47") '  -12.52%  '
Set-TextValue $ws.Range("E48") '  +2.13%  '
Set-TextValue $ws.Range("E49") '  -1.25%  '
Set-TextValue $ws.Range("D50") '2.56'
Set-TextValue $ws.Range("E50") '  -4.39%  '
Set-TextValue $ws.Range("E51") '  +0.12%  '
